$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 74; existing rows 74-81 shift down to 76-83
$ws.Rows("74:75").Insert()

# New row 74: weekly update for "Segunda" quality
$ws.Range("A74").Value = 1
$ws.Range("B74").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C74").Value = "Arica y Parinacota"
$ws.Range("D74").Value = 45218
$ws.Range("E74").Value = 15
$ws.Range("F74").Value = 100112028
$ws.Range("G74").Value = "Sandia"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Segunda"
$ws.Range("J74").Value = 550
$ws.Range("K74").Value = 500
$ws.Range("L74").Value = 550
$ws.Range("M74").Value = 523
$ws.Range("N74").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O74").Value = "Perú"
$ws.Range("P74").Value = 523
$ws.Range("Q74").Value = 1
$ws.Range("R74").Value = "Hortaliza"

# New row 75: weekly update for "Tercera" quality
$ws.Range("A75").Value = 1
$ws.Range("B75").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C75").Value = "Arica y Parinacota"
$ws.Range("D75").Value = 45218
$ws.Range("E75").Value = 15
$ws.Range("F75").Value = 100112028
$ws.Range("G75").Value = "Sandia"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Tercera"
$ws.Range("J75").Value = 350
$ws.Range("K75").Value = 500
$ws.Range("L75").Value = 550
$ws.Range("M75").Value = 529
$ws.Range("N75").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O75").Value = "Perú"
$ws.Range("P75").Value = 529
$ws.Range("Q75").Value = 1
$ws.Range("R75").Value = "Hortaliza"
